$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Introduce new Note values in the same order they were first entered,
# so the workbooks shared-string table indices line up with upstream.
$ws.Range("D16").Value = 'This likely shouldn''t happen; we assume it is due to an expansive definition of "stepchild." Could also be adopted child'
$ws.Range("D58").Value = 'Grandparent'
$ws.Range("D57").Value = 'Aunt/uncle'
$ws.Range("D102").Value = 'Niece/nephew'
$ws.Range("D106").Value = 'Might not be true if half-siblings'
$ws.Range("D120").Value = 'We consider this relationship a grandparent/grandchild relationship, including at birth'
$ws.Range("D99").Value = 'We assume they do not have children with their unmarried partner; if they did, they would be a relative'
$ws.Range("D123").Value = 'Great-grandchild'
$ws.Range("D136").Value = 'Could be parent with expansive "grandchild" definition, but we are hoping that is tracked'
$ws.Range("D140").Value = 'Great-grandparent'
$ws.Range("D152").Value = 'NOT using expansive definition of "grandchild" -- parent-in-law should be direct grandparent of stepchild'
$ws.Range("D26").Value = 'We assume a second-parent adoption'
$ws.Range("D27").Value = 'Possibly biological child (second-parent adoption), but we can''t be sure'

# Remaining cell updates (duplicate Notes + relationship column fixes)
$ws.Range("C26").Value = 'Adopted child'
$ws.Range("C38").Value = 'Adopted child'
$ws.Range("D38").Value = 'We assume a second-parent adoption'
$ws.Range("D39").Value = 'Possibly biological child (second-parent adoption), but we can''t be sure'
$ws.Range("D40").Value = 'This likely shouldn''t happen; we assume it is due to an expansive definition of "stepchild." Could also be adopted child'
$ws.Range("C60").Value = 'Other relative'
$ws.Range("D60").Value = 'Grandparent'
$ws.Range("D73").Value = 'Aunt/uncle'
$ws.Range("D74").Value = 'Grandparent'
$ws.Range("C76").Value = 'Other relative'
$ws.Range("D76").Value = 'Grandparent'
$ws.Range("D86").Value = 'Assuming step-siblings are not reported as "brother or sister"'
$ws.Range("D87").Value = 'Assuming step-siblings are not reported as "brother or sister"'
$ws.Range("D101").Value = 'We assume they do not have children with their unmarried partner; if they did, they would be a relative'
$ws.Range("D103").Value = 'Niece/nephew'
$ws.Range("D115").Value = 'We assume they do not have children with their unmarried partner; if they did, they would be a relative'
$ws.Range("D117").Value = 'We assume they do not have children with their unmarried partner; if they did, they would be a relative'
$ws.Range("C120").Value = 'Grandchild'
$ws.Range("D130").Value = 'Grandparent'
$ws.Range("D131").Value = 'Grandparent'
$ws.Range("D132").Value = 'Grandparent'
$ws.Range("D133").Value = 'Grandparent'
$ws.Range("D135").Value = 'Could be parent, but we are hoping that is tracked'
$ws.Range("D155").Value = 'Great-grandchild'
$ws.Range("D168").Value = 'Could be some kind of spouse/partner, but we need to enforce only 1 spouse/partner'
$ws.Range("D176").Value = 'Assuming foster children are not relatives and "child-in-law" is not used for spouse of a foster child'
$ws.Range("D179").Value = 'We assume they do not have children with their unmarried partner; if they did, they would be a relative'
$ws.Range("D181").Value = 'We assume they do not have children with their unmarried partner; if they did, they would be a relative'
$ws.Range("D224").Value = 'Assuming foster children are not relatives, and not relatives of each other'
$ws.Range("D225").Value = 'Assuming foster children are not relatives, and not relatives of other nonrelatives'

$ws.Range("D26").Select()
